$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.648.19"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.500.04"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.20%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "195.34"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.65%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.627"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("E9").Value = "  -5.91%  "
$ws.Range("E10").Value = "  +0.40%  "
$ws.Range("E11").Value = "  +0.71%  "
$ws.Range("E12").Value = "  -2.14%  "
$ws.Range("E13").Value = "  -0.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.062.05"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.39%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "593.34"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.50%  "
$ws.Range("B16").Value = "Uniswap"
$ws.Range("C16").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "12.78"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.41%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.841.36"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.01"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.75%  "
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.122"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.74%  "
$ws.Range("B20").Value = "WrappedEther"
$ws.Range("C20").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.503.25"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.988"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "18.18"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.63%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.32"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.79%  "
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.65"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "102.14"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.26%  "
$ws.Range("E26").Value = "  +3.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.83"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.91%  "
$ws.Range("E28").Value = "  -1.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.25"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.56%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.30"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.67%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.02"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.39"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.29%  "
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.10"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.52%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0829"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.69%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.712.85"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.08"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.86%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.63"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.07%  "
$ws.Range("E40").Value = "  -1.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.44"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.86%  "
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "471.33"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -9.02%  "
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.133"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.93%  "
$ws.Range("E44").Value = "  -1.90%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.139"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.61%  "
$ws.Range("E46").Value = "  -4.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.28"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.30%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.41"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.16%  "
$ws.Range("E50").Value = "  +1.69%  "
$ws.Range("E51").Value = "  +9.73%  "
